$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: Level (F11) 1 -> 7
$ws.Range("F11").Value = 7

# Row 17: clear Name (B17), Level (F17) 1 -> 99
$ws.Range("B17").ClearContents()
$ws.Range("F17").Value = 99

# Row 20: Level (F20) 5 -> 12
$ws.Range("F20").Value = 12

# Row 21: Level (F21) 10 -> 50
$ws.Range("F21").Value = 50

$ws.Range("F17").Select()
